$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the Requisitos table (rows 20-21) ---
$ws.Range("A20").Value2 = 17
$ws.Range("B20").Value2 = "O usuário precisa se autenticar para utilizar as funcionalidades"
$ws.Range("B20").WrapText = $true
$ws.Range("C20").Value2 = "Autenticar Usuário"
$ws.Range("D20").Value2 = "Usuário"

$ws.Range("A21").Value2 = 18
$ws.Range("B21").Value2 = "Permitir que o usuário filtre os itens que são exibidos na lista de itens de uma viagem"
$ws.Range("B21").WrapText = $true
$ws.Range("C21").Value2 = "Pesquisar Itens"
$ws.Range("D21").Value2 = "Usuário"
$ws.Rows.Item(21).RowHeight = 30

# --- New "Diagrama de Casos de Uso" table (rows 24-50) ---
$ws.Range("A24").Value2 = "UC#"
$ws.Range("B24").Value2 = "NOME DO UC"
$ws.Range("B24").WrapText = $true
$ws.Range("C24").Value2 = "ATOR(ES)"

$ws.Range("A25").Value2 = 1
$ws.Range("B25").Value2 = "Cadastrar Usuário"
$ws.Range("C25").Value2 = "Usuário"

$ws.Range("A26").Value2 = 2
$ws.Range("B26").Value2 = "Criar Viagem"
$ws.Range("C26").Value2 = "Usuário"

$ws.Range("A27").Value2 = 3
$ws.Range("B27").Value2 = "Listar Itens de uma Viagem"
$ws.Range("C27").Value2 = "Usuário"

$ws.Range("A28").Value2 = 4
$ws.Range("B28").Value2 = "Adicionar Voo"
$ws.Range("C28").Value2 = "Usuário"

$ws.Range("A29").Value2 = 5
$ws.Range("B29").Value2 = "Adicionar Hospedagem"
$ws.Range("C29").Value2 = "Usuário"

$ws.Range("A30").Value2 = 6
$ws.Range("B30").Value2 = "Adicionar Trem"
$ws.Range("C30").Value2 = "Usuário"

$ws.Range("A31").Value2 = 7
$ws.Range("B31").Value2 = "Adicionar Atração Turística"
$ws.Range("C31").Value2 = "Usuário"

$ws.Range("A32").Value2 = 8
$ws.Range("B32").Value2 = "Adicionar Nota de Viagem"
$ws.Range("C32").Value2 = "Usuário"

$ws.Range("A33").Value2 = 9
$ws.Range("B33").Value2 = "Gerenciar Voo"
$ws.Range("C33").Value2 = "Usuário"

$ws.Range("A34").Value2 = 10
$ws.Range("B34").Value2 = "Gerenciar Hospedagem"
$ws.Range("B34").WrapText = $true
$ws.Range("C34").Value2 = "Usuário"

$ws.Range("A35").Value2 = 11
$ws.Range("B35").Value2 = "Gerenciar Trem"
$ws.Range("B35").WrapText = $true
$ws.Range("C35").Value2 = "Usuário"

$ws.Range("A36").Value2 = 12
$ws.Range("B36").Value2 = "Gerenciar Atração Turística"
$ws.Range("B36").WrapText = $true
$ws.Range("C36").Value2 = "Usuário"

$ws.Range("A37").Value2 = 13
$ws.Range("B37").Value2 = "Gerenciar Nota de Viagem"
$ws.Range("B37").WrapText = $true
$ws.Range("C37").Value2 = "Usuário"

$ws.Range("A38").Value2 = 14
$ws.Range("B38").Value2 = "Remover Viagem"
$ws.Range("C38").Value2 = "Usuário"

$ws.Range("A39").Value2 = 15
$ws.Range("B39").Value2 = "Enviar E-mail de Ativação de Cadastro"
$ws.Range("B39").WrapText = $true
$ws.Range("C39").Value2 = "Servidor de e-mail"

$ws.Range("A40").Value2 = 16
$ws.Range("B40").Value2 = "Ativar Cadastro"
$ws.Range("C40").Value2 = "Usuário"

$ws.Range("A41").Value2 = 17
$ws.Range("B41").Value2 = "Compartilhar Viagem"
$ws.Range("B41").WrapText = $true
$ws.Range("C41").Value2 = "Usuário"

$ws.Range("A42").Value2 = 18
$ws.Range("B42").Value2 = "Enviar E-mail de Compartilhamento"
$ws.Range("C42").Value2 = "Servidor de e-mail"

$ws.Range("A43").Value2 = 19
$ws.Range("B43").Value2 = "Aceitar Compartilhamento"
$ws.Range("B43").WrapText = $true
$ws.Range("C43").Value2 = "Usuário"

$ws.Range("A44").Value2 = 20
$ws.Range("B44").Value2 = "Listar Viagens Compartilhadas"
$ws.Range("C44").Value2 = "Usuário"

$ws.Range("A45").Value2 = 21
$ws.Range("B45").Value2 = "Listar Viagens Futuras"
$ws.Range("B45").WrapText = $true
$ws.Range("C45").Value2 = "Usuário"

$ws.Range("A46").Value2 = 22
$ws.Range("B46").Value2 = "Listar Histórico de Viagens"
$ws.Range("B46").WrapText = $true
$ws.Range("C46").Value2 = "Usuário"

$ws.Range("A47").Value2 = 23
$ws.Range("B47").Value2 = "Emitir Alerta de Início de Item"
$ws.Range("C47").Value2 = "Sistema"

$ws.Range("A48").Value2 = 24
$ws.Range("B48").Value2 = "Imprimir Itinerário"
$ws.Range("C48").Value2 = "Usuário"

$ws.Range("A49").Value2 = 25
$ws.Range("B49").Value2 = "Autenticar Usuário"
$ws.Range("C49").Value2 = "Usuário"

$ws.Range("A50").Value2 = 26
$ws.Range("B50").Value2 = "Pesquisar Itens"
$ws.Range("C50").Value2 = "Usuário"

# --- Update view/selection state ---
$ws.Activate() | Out-Null
$ws.Range("B55").Select() | Out-Null

Write-Host "Edit complete"
